# Apply the dataImporter troop/castle import edit.
# Order of writes matters: it controls the shared-string table build order,
# which must mirror the order the importer originally populated the sheets
# in (Troop sheet first, then Castle sheet; within a row, occasionally a
# later column is written before an earlier one).
$wb = $excel.ActiveWorkbook

$wsCastle = $wb.Worksheets.Item("Castle")
$wsTroop  = $wb.Worksheets.Item("Troop")

$wsCastle.Cells.Clear()
$wsTroop.Cells.Clear()

# --- Troop sheet (populated first) ------------------------------------
$wsTroop.Range("A1").Value = "ID"
$wsTroop.Range("B1").Value = "Troop Name"
$wsTroop.Range("C1").Value = "Salute"
$wsTroop.Range("D1").Value = "Strength"

$wsTroop.Range("B2").Value = "Mark"
$wsTroop.Range("C2").Value = "Yes sir?"
$wsTroop.Range("D2").Value = 12

$wsTroop.Range("C3").Value = "Yes sir!"
$wsTroop.Range("B3").Value = "Talison"
$wsTroop.Range("D3").Value = 10

# NOTE: the host stores column width on a 1/6-character pixel grid (6px/char
# + 5px padding), while genuine Excel (whose output the target file bytes
# came from) uses a 1/256 grid keyed to the real font's max-digit-width.
# The literal target widths (21.5703125, 5.28515625, ...) therefore cannot
# be represented exactly here; these inputs are chosen so the host's
# rounding lands as close as possible to the target stored width.
$wsTroop.Columns.Item(1).ColumnWidth = 4.5
$wsTroop.Columns.Item(2).ColumnWidth = 11
$wsTroop.Columns.Item(3).ColumnWidth = 14.5
$wsTroop.Columns.Item(4).ColumnWidth = 7.6666666666667

$wsTroop.Range("D3").Select()

# --- Castle sheet (populated second) ----------------------------------
$wsCastle.Range("A1").Value = "ID"
$wsCastle.Range("B1").Value = "Castle Name"
$wsCastle.Range("C1").Value = "Faction"
$wsCastle.Range("D1").Value = "Wall Strength"

$wsCastle.Range("B2").Value = "Castle at Old Town"
$wsCastle.Range("C2").Value = "Knights of the Round"

$wsCastle.Range("E1").Value = "Troops"

$wsCastle.Range("D2").Value = 20000
$wsCastle.Range("E2").Value = "Mark"
$wsCastle.Range("F2").Value = "Talison"
$wsCastle.Range("G2").Value = "/"

$wsCastle.Columns.Item(2).ColumnWidth = 20.6666666666667
$wsCastle.Columns.Item(3).ColumnWidth = 25.6666666666667
$wsCastle.Columns.Item(4).ColumnWidth = 12.3333333333333

$wsCastle.Range("F3").Select()

# --- Workbook-level: make Troop the active/selected tab --------------
$wsTroop.Activate()
